$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 251..354 down to 252..355 (copy from bottom up so we never
# overwrite a source row before it has been copied).
for ($r = 354; $r -ge 251; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# The newly freed row 251 still holds the original row's data (since it was
# copied down, not moved) -- only its Fecha (date, column D) changes.
$ws.Cells.Item(251, 4).Value = 44755
